# Writing-log workbook: add a new day of data (row 5), fill the elapsed-time
# formula in column D down through row 20 (as Excel's fill handle would,
# producing a shared formula), then clear the filled-but-unused cells while
# keeping their number format, and finally update the sheet's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column D down to row 20 with the same "0" (integer) format
# used by the existing D column cells, so the later fill/clear keeps style.
$ws.Range("D5:D20").NumberFormat = "0"

# New row of data (continuing the daily log, 2019-07-18 = serial 43664).
$ws.Range("A5").Value = 43664
$ws.Range("A5").NumberFormat = "d-mmm"
$ws.Range("B5").Formula = "=C4"
$ws.Range("C5").Value = 5445
$ws.Range("E5").Value = 0.5

# Fill the elapsed-time formula from D5 down through D20 (creates the
# shared formula group D5:D20, matching a fill-handle drag in Excel).
$ws.Range("D5:D20").Formula = "=C5-B5"

# The drag went further than needed; clear the contents of the cells
# that shouldn't actually hold a value, keeping their number format.
$ws.Range("D6:D20").ClearContents()

# Leave the selection where the user's cursor ended up.
$ws.Range("F13").Select() | Out-Null
